$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Split the "Colab" run inside the "Google Colab Notebook" hyperlink
#    into two runs ("C" + "olab") that keep identical run formatting.
# ------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Colab", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.Collapse(0)
$rng.Find.Execute("Colab", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$splitPoint = $d.Range($rng.Start, $rng.Start + 1)
$splitPoint.Bold = 1
$splitPoint.Bold = 0

# ------------------------------------------------------------------
# 2) Add the (built-in) "FollowedHyperlink" character style to the
#    style sheet, matching Word's own definition for it, without
#    leaving it applied anywhere in the document body.
# ------------------------------------------------------------------
$marker = "ZZTMP_STYLE_MINT_MARKER_ZZ"
$endRng = $d.Content
$endRng.Collapse(0)
$endRng.InsertAfter($marker)

$markerRng = $d.Content
$markerRng.Find.Execute($marker, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$markerRng.Style = "FollowedHyperlink"

$followedHyperlink = $d.Styles("FollowedHyperlink")
$followedHyperlink.Priority = 99
$followedHyperlink.UnhideWhenUsed = $true
$followedHyperlink.QuickStyle = $false

$markerRng2 = $d.Content
$markerRng2.Find.Execute($marker, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$markerRng2.Delete()

Write-Host "done"
